# Update PEA (column B) projected values for rows 9-24 ("Tasas modelo solo mujeres")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = 1806843
    10 = 1829452
    11 = 1847731
    12 = 1850474
    13 = 1834113
    14 = 1810185
    15 = 1775249
    16 = 1728792
    17 = 1678099
    18 = 1629919
    19 = 1581330
    20 = 1522646
    21 = 1462988
    22 = 1410596
    23 = 1355779
    24 = 1301961
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
